$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 500
$ws.Range("J7").Value = 500
$ws.Range("L7").Value = 500
$ws.Range("N7").Value = -724
$ws.Range("H14").Value = 500
$ws.Range("J14").Value = 500
$ws.Range("L14").Value = 500
$ws.Range("N14").Value = -882
$ws.Range("H113").Value = 75463.31
$ws.Range("I113").Value = 134964.88
$ws.Range("K113").Value = 134964.88
$ws.Range("M113").Value = -131710.88
$ws.Range("H135").Value = 13335182
$ws.Range("I135").Value = 17545292
$ws.Range("J135").Value = 3166.1667
$ws.Range("K135").Value = 157907628
$ws.Range("L135").Value = 28495.5003
$ws.Range("M135").Value = -157905093
$ws.Range("N135").Value = -33565.5003
$ws.Range("H137").Value = 1408733.6
$ws.Range("I137").Value = 51454.727
$ws.Range("J137").Value = 3274992.2
$ws.Range("K137").Value = 154364.181
$ws.Range("L137").Value = 9824976.600000001
$ws.Range("M137").Value = -151814.181
$ws.Range("N137").Value = -9830076.600000001
$ws.Range("H138").Value = 3668.146
$ws.Range("I138").Value = 1396.9166
$ws.Range("J138").Value = 4506.754
$ws.Range("K138").Value = 4190.7498
$ws.Range("L138").Value = 13520.262
$ws.Range("M138").Value = 949.2502000000004
$ws.Range("N138").Value = -23800.262
$ws.Range("H141").Value = 5553
$ws.Range("I141").Value = 5553
$ws.Range("K141").Value = 16659
$ws.Range("M141").Value = -11479

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 797.5714
$ws.Range("I5").Value = 1079.7
$ws.Range("J5").Value = 92.25
$ws.Range("K5").Value = 1079.7
$ws.Range("L5").Value = 92.25
$ws.Range("M5").Value = -967.7
$ws.Range("N5").Value = -316.25
$ws.Range("H45").Value = 43032
$ws.Range("I45").Value = 38794.75
$ws.Range("K45").Value = 38794.75
$ws.Range("M45").Value = -38417.75
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null
$ws.Range("H132").Value = 3067.56
$ws.Range("I132").Value = 2583.0952
$ws.Range("K132").Value = 7749.285600000001
$ws.Range("M132").Value = -5219.285600000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 797.5714
$ws.Range("I4").Value = 1079.7
$ws.Range("J4").Value = 92.25
$ws.Range("K4").Value = 1079.7
$ws.Range("L4").Value = 92.25
$ws.Range("M4").Value = -964.7
$ws.Range("N4").Value = -322.25
$ws.Range("H86").Value = 3106.9412
$ws.Range("I86").Value = 1844.3334
$ws.Range("K86").Value = 1844.3334
$ws.Range("M86").Value = -721.3334
$ws.Range("H89").Value = 3106.9412
$ws.Range("I89").Value = 1844.3334
$ws.Range("K89").Value = 9221.666999999999
$ws.Range("M89").Value = -3605.666999999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2224007
$ws.Range("I31").Value = 2945876.2
$ws.Range("K31").Value = 2945876.2
$ws.Range("M31").Value = -2945581.2
$ws.Range("H34").Value = 2224007
$ws.Range("I34").Value = 2945876.2
$ws.Range("K34").Value = 2945876.2
$ws.Range("M34").Value = -2945674.2
$ws.Range("H58").Value = 4233.5
$ws.Range("I58").Value = 3786.9
$ws.Range("K58").Value = 3786.9
$ws.Range("M58").Value = -3583.9
$ws.Range("H132").Value = 2796.8572
$ws.Range("I132").Value = 2525.7
$ws.Range("K132").Value = 7577.099999999999
$ws.Range("M132").Value = -5047.099999999999
$ws.Range("H136").Value = 4233.5
$ws.Range("I136").Value = 3786.9
$ws.Range("K136").Value = 11360.7
$ws.Range("M136").Value = -8810.700000000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 357.05884
$ws.Range("I113").Value = 283.18182
$ws.Range("J113").Value = 492.5
$ws.Range("K113").Value = 849.54546
$ws.Range("L113").Value = 1477.5
$ws.Range("M113").Value = 1320.45454
$ws.Range("N113").Value = -5817.5
$ws.Range("H114").Value = 1269.1818
$ws.Range("I114").Value = 592.75
$ws.Range("J114").Value = 1655.7142
$ws.Range("K114").Value = 1778.25
$ws.Range("L114").Value = 4967.142599999999
$ws.Range("M114").Value = 1475.75
$ws.Range("N114").Value = -11475.1426

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 20024
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 20024
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 20024
$ws.Range("M38").Value = $null
$ws.Range("N38").Value = -20950
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30518
$ws.Range("H102").Value = 3932.3
$ws.Range("I102").Value = 3924.7778
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 3924.7778
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -2302.7778
$ws.Range("N102").Value = -7244
$ws.Range("H113").Value = 4114.2104
$ws.Range("I113").Value = 3531.3333
$ws.Range("J113").Value = 5113.4287
$ws.Range("K113").Value = 3531.3333
$ws.Range("L113").Value = 5113.4287
$ws.Range("M113").Value = -1361.3333
$ws.Range("N113").Value = -9453.4287

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1539.5714
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4590
$ws.Range("H27").Value = 1539.5714
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("N27").Value = -4214
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31132
$ws.Range("H93").Value = 858373
$ws.Range("I93").Value = 1114315.9
$ws.Range("J93").Value = 5230
$ws.Range("K93").Value = 1114315.9
$ws.Range("L93").Value = 5230
$ws.Range("M93").Value = -1113067.9
$ws.Range("N93").Value = -7726

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 20010
$ws.Range("I43").Value = 15000
$ws.Range("K43").Value = 15000
$ws.Range("M43").Value = -14851
$ws.Range("H122").Value = 5861.7144
$ws.Range("I122").Value = 2466.5
$ws.Range("K122").Value = 7399.5
$ws.Range("M122").Value = -4949.5
$ws.Range("H126").Value = 2968.087
$ws.Range("I126").Value = 2441.238
$ws.Range("K126").Value = 7323.714
$ws.Range("M126").Value = -4853.714
$ws.Range("H132").Value = 8581.629999999999
$ws.Range("I132").Value = 7205.3335
$ws.Range("J132").Value = 11334.223
$ws.Range("K132").Value = 21616.0005
$ws.Range("L132").Value = 34002.669
$ws.Range("M132").Value = -19086.0005
$ws.Range("N132").Value = -39062.669
